$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "TX_GDP" column header in B1 to "GDP"
$ws.Range("B1").Value = "GDP"
